$d = $word.ActiveDocument

# Update the date line (first paragraph, outside the table).
$d.Content.Find.Execute("2024-06-29 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-06-30 Sunday", 2)

# Update the practice-problem table. Cells are addressed positionally
# (row, column) because several old values are duplicated across the
# table (e.g. "49÷2=24, 1" and "97÷9=10, 7" each appear twice with
# different replacements), so a single text search/replace would be
# ambiguous.
$t = $d.Tables.Item(1)

$cellValues = @{
    "1,1" = "41÷6=6, 5"
    "1,2" = "33÷8=4, 1"
    "1,3" = "78÷8=9, 6"
    "1,4" = "77÷7=11, 0"
    "1,5" = "56÷7=8, 0"

    "5,1" = "12÷2=6, 0"
    "5,2" = "22÷5=4, 2"
    "5,3" = "15÷3=5, 0"
    "5,4" = "97÷9=10, 7"
    "5,5" = "70÷2=35, 0"

    "9,1" = "72÷2=36, 0"
    "9,2" = "33÷2=16, 1"
    "9,3" = "10÷9=1, 1"
    "9,4" = "32÷2=16, 0"
    "9,5" = "17÷6=2, 5"

    "13,1" = "31÷6=5, 1"
    "13,2" = "93÷3=31, 0"
    "13,3" = "57÷8=7, 1"
    "13,4" = "57÷9=6, 3"
    "13,5" = "40÷4=10, 0"

    "17,1" = "63÷5=12, 3"
    "17,2" = "86÷4=21, 2"
    "17,3" = "24÷9=2, 6"
    "17,4" = "69÷2=34, 1"
    "17,5" = "53÷2=26, 1"
}

foreach ($key in $cellValues.Keys) {
    $parts = $key.Split(",")
    $row = [int]$parts[0]
    $col = [int]$parts[1]
    $cell = $t.Cell($row, $col)
    $cellRange = $cell.Range
    $cellRange.MoveEnd(1, -1) | Out-Null
    $cellRange.Text = $cellValues[$key]
}
